# Update NATMI LR-pair TPM-derived metrics for Rgmb-Neo1 (OldD0) with new TPM values.
# Only the "Ligand average/total expression" (G/H) and "Receptor average/total
# expression" (M/N) for the ECs cluster changed upstream; every other touched column
# (specificity + edge-weight columns) is a derived value recomputed from those, so we
# write the refreshed values for each affected cell directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 16.57637
$ws.Range("H2").Value = 49.72911
$ws.Range("I2").Value = 0.1853914334114506
$ws.Range("J2").Value = 0.1853914334114506
$ws.Range("M2").Value = 4.331589999999999
$ws.Range("N2").Value = 12.99477
$ws.Range("O2").Value = 0.1478799966101367
$ws.Range("P2").Value = 0.1478799966101367
$ws.Range("Q2").Value = 71.80203852829999
$ws.Range("R2").Value = 646.2183467547
$ws.Range("S2").Value = 0.02741568454443371
$ws.Range("T2").Value = 0.02741568454443371

# Row 3
$ws.Range("G3").Value = 16.57637
$ws.Range("H3").Value = 49.72911
$ws.Range("I3").Value = 0.1853914334114506
$ws.Range("J3").Value = 0.1853914334114506
$ws.Range("O3").Value = 0.5404313285772905
$ws.Range("P3").Value = 0.5404313285772904
$ws.Range("Q3").Value = 262.4024341757867
$ws.Range("R3").Value = 2361.62190758208
$ws.Range("S3").Value = 0.1001913386653985
$ws.Range("T3").Value = 0.1001913386653985

# Row 4
$ws.Range("G4").Value = 16.57637
$ws.Range("H4").Value = 49.72911
$ws.Range("I4").Value = 0.1853914334114506
$ws.Range("J4").Value = 0.1853914334114506
$ws.Range("M4").Value = 9.129751
$ws.Range("O4").Value = 0.3116886748125729
$ws.Range("P4").Value = 0.3116886748125729
$ws.Range("Q4").Value = 151.33813058387
$ws.Range("R4").Value = 1362.04317525483
$ws.Range("S4").Value = 0.0577844102016184
$ws.Range("T4").Value = 0.0577844102016184

# Row 5
$ws.Range("I5").Value = 0.5978024790674488
$ws.Range("J5").Value = 0.5978024790674489
$ws.Range("M5").Value = 4.331589999999999
$ws.Range("N5").Value = 12.99477
$ws.Range("O5").Value = 0.1478799966101367
$ws.Range("P5").Value = 0.1478799966101367
$ws.Range("Q5").Value = 231.5286949589066
$ws.Range("R5").Value = 2083.75825463016
$ws.Range("S5").Value = 0.08840302857802566
$ws.Range("T5").Value = 0.08840302857802568

# Row 6
$ws.Range("I6").Value = 0.5978024790674488
$ws.Range("J6").Value = 0.5978024790674489
$ws.Range("O6").Value = 0.5404313285772905
$ws.Range("P6").Value = 0.5404313285772904
$ws.Range("S6").Value = 0.3230711879892192
$ws.Range("T6").Value = 0.3230711879892192

# Row 7
$ws.Range("I7").Value = 0.5978024790674488
$ws.Range("J7").Value = 0.5978024790674489
$ws.Range("M7").Value = 9.129751
$ws.Range("O7").Value = 0.3116886748125729
$ws.Range("P7").Value = 0.3116886748125729
$ws.Range("Q7").Value = 487.9961709972027
$ws.Range("S7").Value = 0.186328262500204
$ws.Range("T7").Value = 0.186328262500204

# Row 8
$ws.Range("I8").Value = 0.2168060875211005
$ws.Range("J8").Value = 0.2168060875211005
$ws.Range("M8").Value = 4.331589999999999
$ws.Range("N8").Value = 12.99477
$ws.Range("O8").Value = 0.1478799966101367
$ws.Range("P8").Value = 0.1478799966101367
$ws.Range("Q8").Value = 83.96892328250665
$ws.Range("R8").Value = 755.7203095425599
$ws.Range("S8").Value = 0.03206128348767735
$ws.Range("T8").Value = 0.03206128348767735

# Row 9
$ws.Range("I9").Value = 0.2168060875211005
$ws.Range("J9").Value = 0.2168060875211005
$ws.Range("O9").Value = 0.5404313285772905
$ws.Range("P9").Value = 0.5404313285772904
$ws.Range("S9").Value = 0.1171688019226727
$ws.Range("T9").Value = 0.1171688019226726

# Row 10
$ws.Range("I10").Value = 0.2168060875211005
$ws.Range("J10").Value = 0.2168060875211005
$ws.Range("M10").Value = 9.129751
$ws.Range("Q10").Value = 176.9824386212427
$ws.Range("S10").Value = 0.06757600211075052
$ws.Range("T10").Value = 0.06757600211075052
